$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Copy formatting of the first item row (14) down into the new rows (15-17) ---
$ws.Range("A14:K14").Copy()
$ws.Range("A15:K15").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A16:K16").PasteSpecial(-4122)
$ws.Range("A17:K17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Merge the description / part-no cells on the new rows like row 14 ---
$ws.Range("E15:H15").Merge()
$ws.Range("J15:K15").Merge()
$ws.Range("E16:H16").Merge()
$ws.Range("J16:K16").Merge()
$ws.Range("E17:H17").Merge()
$ws.Range("J17:K17").Merge()

# --- Fill in the new purchase-request line items ---
$ws.Range("A15").Value = 2
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = "pcs"
$ws.Range("E15").Value = "item 2"

$ws.Range("A16").Value = 3
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "pcs"
$ws.Range("E16").Value = "item 3"

$ws.Range("A17").Value = 4
$ws.Range("B17").Value = 2
$ws.Range("C17").Value = "pcs"
$ws.Range("E17").Value = "item 4"

# --- Update the Requestor name (was "henne", now "Stephine") ---
$ws.Range("I9").Value = "Stephine"

# --- Update the view: scroll down a bit and leave the selection on I23 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I23").Select()
